# Apply updated cryptocurrency price/volume figures to Sheet1 (columns D and E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (never let Excel coerce it to a number),
# then restore the default "Normal" style so no stray number-format style sticks.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "67.092.46"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "2.630.46"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "591.63"
$ws.Range("E5").Value = "  -2.52%  "
Set-TextValue "D6" "165.97"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -0.04%  "
Set-TextValue "D8" "0.535"
$ws.Range("D9").Value = "2.629.80"
$ws.Range("E9").Value = "  -3.00%  "
Set-TextValue "D10" "0.143"
$ws.Range("E10").Value = "  -1.95%  "
Set-TextValue "D11" "0.160"
$ws.Range("E11").Value = "  +1.26%  "
Set-TextValue "D12" "0.361"
$ws.Range("E12").Value = "  -0.55%  "
Set-TextValue "D13" "5.23"
$ws.Range("E13").Value = "  -1.19%  "
Set-TextValue "D14" "27.55"
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("D15").Value = "3.109.81"
$ws.Range("E15").Value = "  -3.08%  "
Set-TextValue "D16" "0.0000182"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D17").Value = "66.894.74"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "2.622.97"
$ws.Range("E18").Value = "  -3.34%  "
Set-TextValue "D19" "12.14"
$ws.Range("E19").Value = "  +2.72%  "
Set-TextValue "D20" "8.03"
$ws.Range("E20").Value = "  +5.65%  "
Set-TextValue "D21" "359.75"
$ws.Range("E21").Value = "  -2.87%  "
Set-TextValue "D22" "4.34"
$ws.Range("E22").Value = "  -3.01%  "
Set-TextValue "D23" "4.66"
$ws.Range("E23").Value = "  -5.52%  "
Set-TextValue "D24" "10.83"
$ws.Range("E24").Value = "  +8.02%  "
Set-TextValue "D25" "1.95"
$ws.Range("E25").Value = "  -6.05%  "
$ws.Range("E26").Value = "  +0.08%  "
Set-TextValue "D27" "71.02"
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("D28").Value = "2.760.30"
$ws.Range("E28").Value = "  -3.48%  "
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  -2.74%  "
Set-TextValue "D31" "552.69"
$ws.Range("E31").Value = "  -4.88%  "
Set-TextValue "D32" "7.93"
$ws.Range("E32").Value = "  -2.22%  "
Set-TextValue "D33" "1.37"
$ws.Range("E33").Value = "  -3.81%  "
Set-TextValue "D34" "1.90"
$ws.Range("E34").Value = "  -3.48%  "
Set-TextValue "D35" "0.134"
$ws.Range("E35").Value = "  +3.10%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -5.03%  "
Set-TextValue "D38" "157.43"
$ws.Range("E38").Value = "  -2.40%  "
Set-TextValue "D39" "19.18"
$ws.Range("E39").Value = "  -3.39%  "
Set-TextValue "D40" "0.367"
$ws.Range("E40").Value = "  -2.58%  "
Set-TextValue "D41" "5.21"
$ws.Range("E41").Value = "  -3.06%  "
Set-TextValue "D42" "1.80"
$ws.Range("E42").Value = "  -4.11%  "
Set-TextValue "D43" "17.91"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  +0.00%  "
Set-TextValue "D45" "2.47"
$ws.Range("E45").Value = "  -5.30%  "
Set-TextValue "D46" "40.25"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").Value = "0.0₆0299"
$ws.Range("E47").Value = "  -4.16%  "
Set-TextValue "D48" "0.587"
$ws.Range("E48").Value = "  -1.60%  "
Set-TextValue "D49" "151.84"
$ws.Range("E49").Value = "  -1.72%  "
Set-TextValue "D50" "3.82"
$ws.Range("E50").Value = "  -2.11%  "
Set-TextValue "D51" "1.72"
$ws.Range("E51").Value = "  -2.76%  "
